# Append one new data row (row 85) to the sheet, mirroring the data
# pulled from Adafruit IO, matching the existing "all text" layout used
# by the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 85

$ws.Range("A$row").Value = "2024-09-25T18:06:40Z"
$ws.Range("B$row").Value = "temperature"

# Column C holds a numeric-looking value ("25") but the sheet stores every
# cell as plain text, so force text formatting before assigning it to keep
# Excel from auto-converting it into a number. Reset the style back to
# Normal afterwards so no extra formatting is left behind on the cell.
$ws.Range("C$row").NumberFormat = "@"
$ws.Range("C$row").Value = "25"
$ws.Range("C$row").Style = "Normal"

$ws.Range("D$row").Value = "N/A"
$ws.Range("E$row").Value = "N/A"
$ws.Range("F$row").Value = "N/A"
